$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.78343898000433
$ws.Range("C2").Value = 8.320682808848231
$ws.Range("D2").Value = 3.944530293171242
$ws.Range("E2").Value = 11.80144205942174
$ws.Range("F2").Value = 21.68843984460198
$ws.Range("M2").Value = 14.2283810455674
$ws.Range("N2").Value = 16.55318618368922
$ws.Range("O2").Value = 19.17549271906794
$ws.Range("B3").Value = 11.22455376164821
$ws.Range("C3").Value = 7.952968402872831
$ws.Range("D3").Value = 3.908799376978055
$ws.Range("E3").Value = 11.68535303170099
$ws.Range("F3").Value = 21.59894531517007
$ws.Range("M3").Value = 13.95445261415489
$ws.Range("N3").Value = 16.61195351724882
$ws.Range("O3").Value = 19.17683231908321
$ws.Range("B4").Value = 10.86815566977344
$ws.Range("C4").Value = 7.716795311432813
$ws.Range("D4").Value = 3.886432356697452
$ws.Range("E4").Value = 11.61800286406398
$ws.Range("F4").Value = 21.55145764633494
$ws.Range("M4").Value = 13.78690265557993
$ws.Range("N4").Value = 16.64987138642309
$ws.Range("O4").Value = 19.18361076711285
$ws.Range("B5").Value = 10.71978362971639
$ws.Range("C5").Value = 7.618019270520094
$ws.Range("D5").Value = 3.877214260738068
$ws.Range("E5").Value = 11.59157444985732
$ws.Range("F5").Value = 21.53399755154097
$ws.Range("M5").Value = 13.7188835448241
$ws.Range("N5").Value = 16.66578571420689
$ws.Range("O5").Value = 19.18786801577232
$ws.Range("B6").Value = 10.69496330934202
$ws.Range("C6").Value = 7.601467156894174
$ws.Range("D6").Value = 3.875677487036041
$ws.Range("E6").Value = 11.58724828360266
$ws.Range("F6").Value = 21.53121293854106
$ws.Range("M6").Value = 13.70760751819025
$ws.Range("N6").Value = 16.66845624462312
$ws.Range("O6").Value = 19.18866512065584
$ws.Range("B7").Value = 10.86616709788746
$ws.Range("C7").Value = 7.715473326262671
$ws.Range("D7").Value = 3.886308451010485
$ws.Range("E7").Value = 11.61764228522587
$ws.Range("F7").Value = 21.55121449715303
$ws.Range("M7").Value = 13.7859841487943
$ws.Range("N7").Value = 16.65008413844098
$ws.Range("O7").Value = 19.18366213319534
$ws.Range("B8").Value = 11.5935899825638
$ws.Range("C8").Value = 8.196096751481416
$ws.Range("D8").Value = 3.932301681410808
$ws.Range("E8").Value = 11.76061748278038
$ws.Range("F8").Value = 21.6560432843785
$ws.Range("M8").Value = 14.13385486837898
$ws.Range("N8").Value = 16.57306912862203
$ws.Range("O8").Value = 19.17471718070557
$ws.Range("B9").Value = 12.90794741251493
$ws.Range("C9").Value = 9.053186629045424
$ws.Range("D9").Value = 4.018912906889112
$ws.Range("E9").Value = 12.07075770053318
$ws.Range("F9").Value = 21.92001595147825
$ws.Range("M9").Value = 14.81694917484901
$ws.Range("N9").Value = 16.43654429983514
$ws.Range("O9").Value = 19.20450340591557
$ws.Range("B10").Value = 13.79749841492937
$ws.Range("C10").Value = 9.627619568836602
$ws.Range("D10").Value = 4.080121807558134
$ws.Range("E10").Value = 12.31473336545947
$ws.Range("F10").Value = 22.14827884225035
$ws.Range("M10").Value = 15.31395588563286
$ws.Range("N10").Value = 16.34500352820977
$ws.Range("O10").Value = 19.25526019056306
$ws.Range("B11").Value = 14.18442696672209
$ws.Range("C11").Value = 9.876441002296151
$ws.Range("D11").Value = 4.107391319341943
$ws.Range("E11").Value = 12.42875875494141
$ws.Range("F11").Value = 22.25924311960859
$ws.Range("M11").Value = 15.53791975109943
$ws.Range("N11").Value = 16.30524626829387
$ws.Range("O11").Value = 19.28459914942795
$ws.Range("B12").Value = 14.32831695990071
$ws.Range("C12").Value = 9.9688344007156
$ws.Range("D12").Value = 4.117630846579344
$ws.Range("E12").Value = 12.47233385142225
$ws.Range("F12").Value = 22.30225552077454
$ws.Range("M12").Value = 15.62233863622219
$ws.Range("N12").Value = 16.29046114303005
$ws.Range("O12").Value = 19.29660382663606
$ws.Range("B13").Value = 14.29744585841646
$ws.Range("C13").Value = 9.949017686367748
$ws.Range("D13").Value = 4.115429509389172
$ws.Range("E13").Value = 12.46293221333019
$ws.Range("F13").Value = 22.29294841803572
$ws.Range("M13").Value = 15.60417623693131
$ws.Range("N13").Value = 16.29363338762601
$ws.Range("O13").Value = 19.29397869214981
$ws.Range("B14").Value = 14.19631803800133
$ws.Range("C14").Value = 9.884079107880453
$ws.Range("D14").Value = 4.108235495598915
$ws.Range("E14").Value = 12.43233600130041
$ws.Range("F14").Value = 22.26276207378093
$ws.Range("M14").Value = 15.54487326041681
$ws.Range("N14").Value = 16.30402447978973
$ws.Range("O14").Value = 19.28556887881007
$ws.Range("B15").Value = 14.13402945837208
$ws.Range("C15").Value = 9.844063152630817
$ws.Range("D15").Value = 4.103817521773703
$ws.Range("E15").Value = 12.41364530427746
$ws.Range("F15").Value = 22.24440034516854
$ws.Range("M15").Value = 15.50849504872006
$ws.Range("N15").Value = 16.31042446803569
$ws.Range("O15").Value = 19.28053399446596
$ws.Range("B16").Value = 13.77184759915857
$ws.Range("C16").Value = 9.611104132074489
$ws.Range("D16").Value = 4.078327754471063
$ws.Range("E16").Value = 12.30733921804147
$ws.Range("F16").Value = 22.14116777241279
$ws.Range("M16").Value = 15.29926910453794
$ws.Range("N16").Value = 16.34763960658745
$ws.Range("O16").Value = 19.25346826139271
$ws.Range("B17").Value = 13.54505633832474
$ws.Range("C17").Value = 9.464966124799307
$ws.Range("D17").Value = 4.062540385874831
$ws.Range("E17").Value = 12.24287348476386
$ws.Range("F17").Value = 22.07964104083311
$ws.Range("M17").Value = 15.17030869359965
$ws.Range("N17").Value = 16.37095199654325
$ws.Range("O17").Value = 19.23846249190675
$ws.Range("B18").Value = 13.41294679643635
$ws.Range("C18").Value = 9.379737287584121
$ws.Range("D18").Value = 4.053406058119783
$ws.Range("E18").Value = 12.20608310727699
$ws.Range("F18").Value = 22.04492533293594
$ws.Range("M18").Value = 15.09593870589636
$ws.Range("N18").Value = 16.38453815840131
$ws.Range("O18").Value = 19.23042011826202
$ws.Range("B19").Value = 13.36793349707947
$ws.Range("C19").Value = 9.350679739335956
$ws.Range("D19").Value = 4.050304212440578
$ws.Range("E19").Value = 12.19367732295297
$ws.Range("F19").Value = 22.03328769283274
$ws.Range("M19").Value = 15.07072743199312
$ws.Range("N19").Value = 16.38916871425651
$ws.Range("O19").Value = 19.22779829391802
$ws.Range("B20").Value = 13.56937163383238
$ws.Range("C20").Value = 9.480644545601823
$ws.Range("D20").Value = 4.064226582327325
$ws.Range("E20").Value = 12.24970639696943
$ws.Range("F20").Value = 22.08612125535524
$ws.Range("M20").Value = 15.18405760373361
$ws.Range("N20").Value = 16.3684519903201
$ws.Range("O20").Value = 19.23999899416583
$ws.Range("B21").Value = 14.22609372670881
$ws.Range("C21").Value = 9.903203043570166
$ws.Range("D21").Value = 4.110350941734878
$ws.Range("E21").Value = 12.44131242794269
$ws.Range("F21").Value = 22.27160186204868
$ws.Range("M21").Value = 15.56230326213483
$ws.Range("N21").Value = 16.3009650418359
$ws.Range("O21").Value = 19.28801480254889
$ws.Range("B22").Value = 14.63993320616977
$ws.Range("C22").Value = 10.16869196476074
$ws.Range("D22").Value = 4.13998780514592
$ws.Range("E22").Value = 12.56882783048829
$ws.Range("F22").Value = 22.3985935762497
$ws.Range("M22").Value = 15.80719114238375
$ws.Range("N22").Value = 16.25843224671668
$ws.Range("O22").Value = 19.324607862427
$ws.Range("B23").Value = 14.42048868112881
$ws.Range("C23").Value = 10.0279821343331
$ws.Range("D23").Value = 4.124217902045417
$ws.Range("E23").Value = 12.50057459984168
$ws.Range("F23").Value = 22.3302989878423
$ws.Range("M23").Value = 15.67672895778955
$ws.Range("N23").Value = 16.28098911523052
$ws.Range("O23").Value = 19.30460221587694
$ws.Range("B24").Value = 13.55838405297596
$ws.Range("C24").Value = 9.473560104583816
$ws.Range("D24").Value = 4.063464433109481
$ws.Range("E24").Value = 12.24661638753293
$ws.Range("F24").Value = 22.0831895019352
$ws.Range("M24").Value = 15.17784242795896
$ws.Range("N24").Value = 16.36958167125512
$ws.Range("O24").Value = 19.2393025200329
$ws.Range("B25").Value = 12.56526752789005
$ws.Range("C25").Value = 8.830831317804964
$ws.Range("D25").Value = 3.995890858450329
$ws.Range("E25").Value = 11.98387268917402
$ws.Range("F25").Value = 21.84247297828361
$ws.Range("M25").Value = 14.63262776554357
$ws.Range("N25").Value = 16.47193339399939
$ws.Range("O25").Value = 19.19137171584725
